$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" (column E) was previously listed in descending order
# (2104 .. 2004). It is now listed in ascending order (2004 .. 2104).
$periodos = @("2004","2005","2006","2007","2008","2009","2010","2012","2101","2102","2103","2104")

# "Valor Mora" (column F) values follow the same re-ordering as the
# periods above (the block of rows was effectively reversed).
$valores = @(36000,36000,36000,36000,36000,36000,36000,35600,35600,35600,35600,30854)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
